$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseSequence")

# Collapse the old FromLine/ToLine/LineNumbers columns (F,G,H) into a single
# "Lines" column. Deleting the F:G entire columns shifts H (with its
# LineNumbers header/value and its 12.6640625-wide column formatting) left
# into F, preserving the column width exactly.
$ws.Range("F1:G2").EntireColumn.Delete() | Out-Null

# Rename the header and replace the value with the new "1-2" range format.
$ws.Range("F1").Value = "Lines"
$ws.Range("F2").Value = "'1-2"
$ws.Range("F2").NumberFormat = "@"

# Make TestCaseSequence the active sheet/selection, as in the target workbook.
[void]$ws.Activate()
$ws.Range("F2").Select() | Out-Null
